$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.019.18"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -2.52%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.667.72"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "216.93"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.48%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5108"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "
$ws.Range("E7").Value = "  -0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2653"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06406"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.58%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "21.82"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.84%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07416"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.79%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.689.33"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.57%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.500"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.70%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.5836"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.42%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.000008570"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "64.36"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.92%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "26.085.86"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.37%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "4.942"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("E19").Value = "  +0.06%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "10.77"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -2.26%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "190.52"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.73%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.227"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.59%  "
$ws.Range("E23").Value = "  +0.17%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "145.12"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.35%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "7.622"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.32%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.1199"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +3.22%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "15.63"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.69%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.06571"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +16.10%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.328"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.92%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.317"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -2.05%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.540"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.75%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.524"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.62%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.645"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.018"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.6095"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("E36").Value = "  +0.45%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.709"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.80%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "6.264"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +7.05%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01603"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.80%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.084.82"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.8590"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +0.55%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "100.23"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.22%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "1.817.20"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -2.09%  "
$ws.Range("E45").Value = "  +2.51%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "56.32"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.84%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.10%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "8.053"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.55%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.05237"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("E51").Value = "  +3.69%  "
